# framework can now read multiple pdf's in a single execution flow.
# and updated data for klic poc
#
# Updates the "Test Cases" sheet:
#  - drops the old "Groups" (C) and "Execute" (D) columns for the existing
#    101/1898..1973 rows (only D2 keeps a value, now "TestCaseNumber=1974")
#  - appends new TestCaseNumber/Priority rows (1974-2008 / 77-111)
#  - moves the sheet's selection

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

# Drop the Groups/Execute columns for the pre-existing rows (2-83).
$ws.Range("C2:D83").ClearContents()

# D2 now carries the marker that used to live under "Groups"/"Execute".
$d2 = $ws.Cells.Item(2, 4)
$d2.NumberFormat = "@"
$d2.HorizontalAlignment = -4131
$d2.Value = "TestCaseNumber=1974"

# New TestCaseNumber / Priority rows appended below the existing data.
$newRows = @(
    @(1974,77),
    @(1975,78),
    @(1976,79),
    @(1977,80),
    @(1978,81),
    @(1979,82),
    @(1980,83),
    @(1981,84),
    @(1982,85),
    @(1983,86),
    @(1984,87),
    @(1985,88),
    @(1986,89),
    @(1987,90),
    @(1988,91),
    @(1989,92),
    @(1990,93),
    @(1991,94),
    @(1992,95),
    @(1993,96),
    @(1994,97),
    @(1995,98),
    @(1996,99),
    @(1997,100),
    @(1998,101),
    @(1999,102),
    @(2000,103),
    @(2001,104),
    @(2002,105),
    @(2003,106),
    @(2004,107),
    @(2005,108),
    @(2006,109),
    @(2007,110),
    @(2008,111)
)

$row = 84
foreach ($pair in $newRows) {
    $aCell = $ws.Cells.Item($row, 1)
    $aCell.NumberFormat = "@"
    $aCell.HorizontalAlignment = -4131
    $aCell.Value = [string]$pair[0]

    $bCell = $ws.Cells.Item($row, 2)
    $bCell.NumberFormat = "@"
    $bCell.HorizontalAlignment = -4131
    $bCell.Value = [string]$pair[1]

    $row = $row + 1
}

# Move the active selection on the "Test Cases" tab.
$ws.Activate()
$ws.Range("D12").Select()

Write-Host "done"
